{"js": "const replacements = [\n  [\"460\u00f79=51, 1\", \"617\u00f76=102, 5\"],\n  [\"243\u00f72=121, 1\", \"642\u00f79=71, 3\"],\n  [\"516\u00f73=172, 0\", \"669\u00f78=83, 5\"],\n  [\"152\u00f75=30, 2\", \"439\u00f75=87, 4\"],\n  [\"546\u00f72=273, 0\", \"777\u00f73=259, 0\"],\n  [\"681\u00f74=170, 1\", \"485\u00f77=69, 2\"],\n  [\"525\u00f75=105, 0\", \"389\u00f73=129, 2\"],\n  [\"473\u00f75=94, 3\", \"540\u00f78=67, 4\"],\n  [\"237\u00f72=118, 1\", \"265\u00f76=44, 1\"],\n  [\"256\u00f73=85, 1\", \"778\u00f75=155, 3\"],\n  [\"956\u00f79=106, 2\", \"722\u00f79=80, 2\"],\n  [\"409\u00f75=81, 4\", \"689\u00f79=76, 5\"],\n  [\"170\u00f73=56, 2\", \"621\u00f79=69, 0\"],\n  [\"444\u00f73=148, 0\", \"793\u00f79=88, 1\"],\n  [\"696\u00f77=99, 3\", \"281\u00f74=70, 1\"],\n  [\"218\u00f76=36, 2\", \"838\u00f73=279, 1\"],\n  [\"717\u00f77=102, 3\", \"971\u00f77=138, 5\"],\n  [\"670\u00f77=95, 5\", \"247\u00f75=49, 2\"],\n  [\"638\u00f78=79, 6\", \"419\u00f75=83, 4\"],\n  [\"146\u00f74=36, 2\", \"842\u00f79=93, 5\"],\n  [\"721\u00f76=120, 1\", \"766\u00f73=255, 1\"],\n  [\"623\u00f79=69, 2\", \"497\u00f72=248, 1\"],\n  [\"278\u00f72=139, 0\", \"465\u00f78=58, 1\"],\n  [\"433\u00f78=54, 1\", \"702\u00f77=100, 2\"],\n  [\"444\u00f74=111, 0\", \"902\u00f72=451, 0\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$map = @{\n  \"460\u00f79=51, 1\" = \"617\u00f76=102, 5\"\n  \"243\u00f72=121, 1\" = \"642\u00f79=71, 3\"\n  \"516\u00f73=172, 0\" = \"669\u00f78=83, 5\"\n  \"152\u00f75=30, 2\" = \"439\u00f75=87, 4\"\n  \"546\u00f72=273, 0\" = \"777\u00f73=259, 0\"\n  \"681\u00f74=170, 1\" = \"485\u00f77=69, 2\"\n  \"525\u00f75=105, 0\" = \"389\u00f73=129, 2\"\n  \"473\u00f75=94, 3\" = \"540\u00f78=67, 4\"\n  \"237\u00f72=118, 1\" = \"265\u00f76=44, 1\"\n  \"256\u00f73=85, 1\" = \"778\u00f75=155, 3\"\n  \"956\u00f79=106, 2\" = \"722\u00f79=80, 2\"\n  \"409\u00f75=81, 4\" = \"689\u00f79=76, 5\"\n  \"170\u00f73=56, 2\" = \"621\u00f79=69, 0\"\n  \"444\u00f73=148, 0\" = \"793\u00f79=88, 1\"\n  \"696\u00f77=99, 3\" = \"281\u00f74=70, 1\"\n  \"218\u00f76=36, 2\" = \"838\u00f73=279, 1\"\n  \"717\u00f77=102, 3\" = \"971\u00f77=138, 5\"\n  \"670\u00f77=95, 5\" = \"247\u00f75=49, 2\"\n  \"638\u00f78=79, 6\" = \"419\u00f75=83, 4\"\n  \"146\u00f74=36, 2\" = \"842\u00f79=93, 5\"\n  \"721\u00f76=120, 1\" = \"766\u00f73=255, 1\"\n  \"623\u00f79=69, 2\" = \"497\u00f72=248, 1\"\n  \"278\u00f72=139, 0\" = \"465\u00f78=58, 1\"\n  \"433\u00f78=54, 1\" = \"702\u00f77=100, 2\"\n  \"444\u00f74=111, 0\" = \"902\u00f72=451, 0\"\n}\n\nforeach ($table in $d.Tables) {\n  foreach ($row in $table.Rows) {\n    foreach ($cell in $row.Cells) {\n      $r = $cell.Range\n      $txt = $r.Text\n      $txt = $txt.TrimEnd([char]13, [char]7)\n      if ($map.ContainsKey($txt)) {\n        $r.Find.Execute($txt, $false, $false, $false, $false, $false, $true, 1, $false, $map[$txt], 2)\n      }\n    }\n  }\n}\n"}
